$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.290.33'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.909.17'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.11'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5309'
$ws.Range("E7").Value = '  +1.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3825'
$ws.Range("E8").Value = '  +1.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07299'
$ws.Range("E9").Value = '  +0.33%  '

$ws.Range("E10").Value = '  +4.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9023'
$ws.Range("E11").Value = '  -0.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08192'
$ws.Range("E12").Value = '  -0.69%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.97'
$ws.Range("E13").Value = '  -0.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.351'
$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008659'
$ws.Range("E16").Value = '  -0.19%  '

$ws.Range("E17").Value = '  +1.73%  '

$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.266.86'
$ws.Range("E19").Value = '  -33.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.326.46'
$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.061'
$ws.Range("E21").Value = '  -0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.83'
$ws.Range("E22").Value = '  +1.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.520'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.96'
$ws.Range("E24").Value = '  +1.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("E25").Value = '  -1.52%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.02'
$ws.Range("E28").Value = '  +1.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.828'
$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.820'
$ws.Range("E30").Value = '  -1.69%  '

$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8348'
$ws.Range("E32").Value = '  +4.30%  '

$ws.Range("E33").Value = '  -0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.226'
$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.002'
$ws.Range("E35").Value = '  +1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.362'
$ws.Range("E36").Value = '  -1.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.690'
$ws.Range("E37").Value = '  +3.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5749'
$ws.Range("E38").Value = '  +0.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02010'
$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.076'
$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.347'
$ws.Range("E41").Value = '  +3.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.565'
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '117.10'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1526'
$ws.Range("E44").Value = '  +0.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4940'
$ws.Range("E45").Value = '  +1.63%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.08%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.12'
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.640'
$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.77'
$ws.Range("E49").Value = '  +2.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06166'
$ws.Range("E50").Value = '  +3.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.68'
$ws.Range("E51").Value = '  -0.46%  '
